# Add a new worksheet "ODI Bowling Extra" (scraped bowling stats not covered
# by the existing "ODI Bowling" sheet: maiden overs + percent of all wickets)
# as the last sheet in the workbook, matching the "ODI Batting Extra" sheet
# style/layout already used for the equivalent batting data.

$wb = $excel.ActiveWorkbook

$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)

$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "ODI Bowling Extra"

# Data rows: MATCH_CODE, MAIDEN_OVERS, PERCENT_WICKETS_OF_ALL.
# Kept as plain text (not numbers/percentages) to match the rest of the
# workbook's scraped-data convention, so pre-format the whole used range as
# Text before any values are written.
$data = @(
    @("4353", $null, $null),
    @("4359", $null, $null),
    @("4360", "0", "10.00%"),
    @("4362", $null, $null),
    @("4385", $null, $null),
    @("4387", "0", "20.00%"),
    @("4388", "0", "10.00%"),
    @("4398", "0", $null),
    @("4399", "0", "20.00%"),
    @("4400", "1", "20.00%"),
    @("4402", $null, $null),
    @("4406", $null, $null),
    @("4410", "0", "10.00%"),
    @("4435", $null, $null),
    @("4436", $null, $null),
    @("4437", "0", "10.00%"),
    @("4613", "0", $null),
    @("4618", "0", "10.00%"),
    @("4725", "0", "20.00%"),
    @("4732", "0", $null)
)

$lastRow = 1 + $data.Count
$usedRange = $ws.Range("A1:C$lastRow")
$usedRange.NumberFormat = "@"

$ws.Cells.Item(1, 1).Value = "MATCH_CODE"
$ws.Cells.Item(1, 2).Value = "MAIDEN_OVERS"
$ws.Cells.Item(1, 3).Value = "PERCENT_WICKETS_OF_ALL"
$ws.Range("A1:C1").Font.Bold = $true

$rowIdx = 2
foreach ($row in $data) {
    $ws.Cells.Item($rowIdx, 1).Value = $row[0]
    $ws.Cells.Item($rowIdx, 2).Value = $row[1]
    $ws.Cells.Item($rowIdx, 3).Value = $row[2]
    $rowIdx++
}

# Restore focus to the first sheet (the edit didn't change the active tab).
$wb.Worksheets.Item(1).Activate()

Write-Host "Added 'ODI Bowling Extra' sheet with $($data.Count) data rows"
